$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 248.07
$ws.Cells.Item(15, 9).Value = 248.07
$ws.Cells.Item(15, 11).Value = 744.21
$ws.Cells.Item(15, 13).Value = -575.21

$ws.Cells.Item(129, 8).Value = 1072.2118
$ws.Cells.Item(129, 9).Value = 625
$ws.Cells.Item(129, 10).Value = 1094.2963
$ws.Cells.Item(129, 11).Value = 1875
$ws.Cells.Item(129, 12).Value = 3282.8889
$ws.Cells.Item(129, 13).Value = 3125
$ws.Cells.Item(129, 14).Value = -13282.8889

$ws.Cells.Item(132, 8).Value = 1220.127
$ws.Cells.Item(132, 9).Value = 1139.3103
$ws.Cells.Item(132, 10).Value = 2157.6
$ws.Cells.Item(132, 11).Value = 3417.9309
$ws.Cells.Item(132, 12).Value = 6472.799999999999
$ws.Cells.Item(132, 13).Value = -887.9309000000003
$ws.Cells.Item(132, 14).Value = -11532.8

$ws.Cells.Item(135, 8).Value = 2566.322
$ws.Cells.Item(135, 9).Value = 2161.05
$ws.Cells.Item(135, 10).Value = 3419.5264
$ws.Cells.Item(135, 11).Value = 19449.45
$ws.Cells.Item(135, 12).Value = 30775.7376
$ws.Cells.Item(135, 13).Value = -16914.45
$ws.Cells.Item(135, 14).Value = -35845.7376

$ws.Cells.Item(138, 8).Value = 2111.647
$ws.Cells.Item(138, 9).Value = 1280.2069
$ws.Cells.Item(138, 10).Value = 2542.2144
$ws.Cells.Item(138, 11).Value = 3840.620699999999
$ws.Cells.Item(138, 12).Value = 7626.6432
$ws.Cells.Item(138, 13).Value = 1299.379300000001
$ws.Cells.Item(138, 14).Value = -17906.6432

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1558.1555
$ws.Cells.Item(74, 9).Value = 1182.091
$ws.Cells.Item(74, 10).Value = 2592.3333
$ws.Cells.Item(74, 11).Value = 1182.091
$ws.Cells.Item(74, 12).Value = 2592.3333
$ws.Cells.Item(74, 13).Value = -308.0909999999999
$ws.Cells.Item(74, 14).Value = -4340.3333

$ws.Cells.Item(77, 8).Value = 1558.1555
$ws.Cells.Item(77, 9).Value = 1182.091
$ws.Cells.Item(77, 10).Value = 2592.3333
$ws.Cells.Item(77, 11).Value = 5910.455
$ws.Cells.Item(77, 12).Value = 12961.6665
$ws.Cells.Item(77, 13).Value = -1542.455
$ws.Cells.Item(77, 14).Value = -21697.6665

$ws.Cells.Item(132, 8).Value = 2043757.6
$ws.Cells.Item(132, 9).Value = 1936.1052
$ws.Cells.Item(132, 10).Value = 9097323
$ws.Cells.Item(132, 11).Value = 5808.3156
$ws.Cells.Item(132, 12).Value = 27291969
$ws.Cells.Item(132, 13).Value = -3278.3156
$ws.Cells.Item(132, 14).Value = -27297029

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value = 18585.6
$ws.Cells.Item(82, 9).Value = 4752.3335
$ws.Cells.Item(82, 10).Value = 24514.143
$ws.Cells.Item(82, 11).Value = 4752.3335
$ws.Cells.Item(82, 12).Value = 24514.143
$ws.Cells.Item(82, 13).Value = -4369.3335
$ws.Cells.Item(82, 14).Value = -25280.143

$ws.Cells.Item(85, 8).Value = 18585.6
$ws.Cells.Item(85, 9).Value = 4752.3335
$ws.Cells.Item(85, 10).Value = 24514.143
$ws.Cells.Item(85, 11).Value = 4752.3335
$ws.Cells.Item(85, 12).Value = 24514.143
$ws.Cells.Item(85, 13).Value = -3426.3335
$ws.Cells.Item(85, 14).Value = -27166.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 182739.25
$ws.Cells.Item(31, 9).Value = 1660.1666
$ws.Cells.Item(31, 10).Value = 488310.2
$ws.Cells.Item(31, 11).Value = 1660.1666
$ws.Cells.Item(31, 12).Value = 488310.2
$ws.Cells.Item(31, 13).Value = -1365.1666
$ws.Cells.Item(31, 14).Value = -488900.2

$ws.Cells.Item(34, 8).Value = 182739.25
$ws.Cells.Item(34, 9).Value = 1660.1666
$ws.Cells.Item(34, 10).Value = 488310.2
$ws.Cells.Item(34, 11).Value = 1660.1666
$ws.Cells.Item(34, 12).Value = 488310.2
$ws.Cells.Item(34, 13).Value = -1458.1666
$ws.Cells.Item(34, 14).Value = -488714.2

$ws.Cells.Item(134, 8).Value = 224866.89
$ws.Cells.Item(134, 9).Value = 2587.05
$ws.Cells.Item(134, 11).Value = 7761.150000000001
$ws.Cells.Item(134, 13).Value = -5226.150000000001

$ws.Cells.Item(141, 8).Value = 375824.66
$ws.Cells.Item(141, 10).Value = 417826.38
$ws.Cells.Item(141, 12).Value = 417826.38
$ws.Cells.Item(141, 14).Value = -428186.38

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 2381841
$ws.Cells.Item(131, 9).Value = 5263706
$ws.Cells.Item(131, 10).Value = 1170.3043
$ws.Cells.Item(131, 11).Value = 15791118
$ws.Cells.Item(131, 12).Value = 3510.9129
$ws.Cells.Item(131, 13).Value = -15786078
$ws.Cells.Item(131, 14).Value = -13590.9129

$ws.Cells.Item(134, 8).Value = 7949.355
$ws.Cells.Item(134, 9).Value = 6654.2104
$ws.Cells.Item(134, 10).Value = 10000
$ws.Cells.Item(134, 11).Value = 19962.6312
$ws.Cells.Item(134, 12).Value = 30000
$ws.Cells.Item(134, 13).Value = -14892.6312
$ws.Cells.Item(134, 14).Value = -40140

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 5799.4644
$ws.Cells.Item(80, 9).Value = 7693.6113
$ws.Cells.Item(80, 10).Value = 2390
$ws.Cells.Item(80, 11).Value = 7693.6113
$ws.Cells.Item(80, 12).Value = 2390
$ws.Cells.Item(80, 13).Value = -6695.6113
$ws.Cells.Item(80, 14).Value = -4386

$ws.Cells.Item(83, 8).Value = 5799.4644
$ws.Cells.Item(83, 9).Value = 7693.6113
$ws.Cells.Item(83, 10).Value = 2390
$ws.Cells.Item(83, 11).Value = 38468.0565
$ws.Cells.Item(83, 12).Value = 11950
$ws.Cells.Item(83, 13).Value = -33476.0565
$ws.Cells.Item(83, 14).Value = -21934

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2408.6667
$ws.Cells.Item(7, 9).Value = 2166.5557
$ws.Cells.Item(7, 10).Value = 3135
$ws.Cells.Item(7, 11).Value = 2166.5557
$ws.Cells.Item(7, 12).Value = 3135
$ws.Cells.Item(7, 13).Value = -2054.5557
$ws.Cells.Item(7, 14).Value = -3359

$ws.Cells.Item(82, 8).Value = 557946.6
$ws.Cells.Item(82, 9).Value = 910884.5600000001
$ws.Cells.Item(82, 10).Value = 126578.11
$ws.Cells.Item(82, 11).Value = 910884.5600000001
$ws.Cells.Item(82, 12).Value = 126578.11
$ws.Cells.Item(82, 13).Value = -910523.5600000001
$ws.Cells.Item(82, 14).Value = -127300.11

$ws.Cells.Item(85, 8).Value = 557946.6
$ws.Cells.Item(85, 9).Value = 910884.5600000001
$ws.Cells.Item(85, 10).Value = 126578.11
$ws.Cells.Item(85, 11).Value = 910884.5600000001
$ws.Cells.Item(85, 12).Value = 126578.11
$ws.Cells.Item(85, 13).Value = -909636.5600000001
$ws.Cells.Item(85, 14).Value = -129074.11

$ws.Cells.Item(126, 8).Value = 2408.6667
$ws.Cells.Item(126, 9).Value = 2166.5557
$ws.Cells.Item(126, 10).Value = 3135
$ws.Cells.Item(126, 11).Value = 6499.6671
$ws.Cells.Item(126, 12).Value = 9405
$ws.Cells.Item(126, 13).Value = -4029.6671
$ws.Cells.Item(126, 14).Value = -14345

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1634.5
$ws.Cells.Item(122, 9).Value = 1276.3
$ws.Cells.Item(122, 10).Value = 2146.2144
$ws.Cells.Item(122, 11).Value = 3828.9
$ws.Cells.Item(122, 12).Value = 6438.6432
$ws.Cells.Item(122, 13).Value = -1378.9
$ws.Cells.Item(122, 14).Value = -11338.6432

$ws.Cells.Item(136, 8).Value = 1925.8088
$ws.Cells.Item(136, 9).Value = 1808.6666
$ws.Cells.Item(136, 10).Value = 2155
$ws.Cells.Item(136, 11).Value = 5425.9998
$ws.Cells.Item(136, 12).Value = 6465
$ws.Cells.Item(136, 13).Value = -2875.9998
$ws.Cells.Item(136, 14).Value = -11565

